# road_inspection_database.xlsx - "New Changes with the old project"
#
# 1. Rename the sheet from Road_Inspection_Data -> Sheet1
# 2. Simplify the header styling: drop the dark-blue fill and the bold
#    white font, leaving a plain bold header (border + centered alignment
#    are kept).
# 3. Refresh the two detections that survive (rows 2 & 3) with new data,
#    and drop the two extra detection rows (old rows 4 & 5) so the sheet
#    only spans A1:Q3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the worksheet -------------------------------------------------
$ws.Name = "Sheet1"

# --- 2. Simplify the header row formatting -----------------------------------
$header = $ws.Range("A1:Q1")
$header.Font.Color = 0           # was white (FFFFFF) -> plain/black bold text
$header.Interior.Pattern = -4142 # xlNone: drop the dark "002c3e50" fill

# --- 3. Drop the last two detections (old rows 4 and 5) -----------------------
$ws.Rows("4:5").Delete()

# --- 4. Update row 2 (first detection) ----------------------------------------
$ws.Range("A2").Value = "DET_20250927_080330_1"
$ws.Range("B2").Value = "2025-09-27 08:03:30"
$ws.Range("C2").Value = "Tumakuru"
$ws.Range("D2").Value = "Mg road"
$ws.Range("H2").Value = 0.819
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = "Jnaneshwari P"

# --- 5. Update row 3 (second detection) ---------------------------------------
$ws.Range("A3").Value = "DET_20250927_080330_2"
$ws.Range("B3").Value = "2025-09-27 08:03:30"
$ws.Range("C3").Value = "Tumakuru"
$ws.Range("D3").Value = "Mg road"
$ws.Range("G3").Value = "pothole"
$ws.Range("H3").Value = 0.287
$ws.Range("K3").Value = "Hot Mix Asphalt Patching"
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = "Jnaneshwari P"
$ws.Range("P3").Value = 17
